# Auto-generated Excel COM-interop script to apply Ultima_Profits market-data updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1800.125
$ws.Range("I40").Value = 1771.5714
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1771.5714
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1596.5714
$ws.Range("N40").Value = -2350
# Row 68
$ws.Range("H68").Value = 15000
$ws.Range("J68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16498
# Row 71
$ws.Range("H71").Value = 15000
$ws.Range("J71").Value = 15000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -52488
# Row 116
$ws.Range("H116").Value = 2249.375
$ws.Range("I116").Value = 1801.25
$ws.Range("J116").Value = 2697.5
$ws.Range("K116").Value = 1801.25
$ws.Range("L116").Value = 2697.5
$ws.Range("M116").Value = 1640.75
$ws.Range("N116").Value = -9581.5
# Row 125
$ws.Range("H125").Value = 3629.8572
$ws.Range("I125").Value = 4501.8
$ws.Range("J125").Value = 1450
$ws.Range("K125").Value = 40516.2
$ws.Range("L125").Value = 13050
$ws.Range("M125").Value = -38056.2
$ws.Range("N125").Value = -17970
# Row 129
$ws.Range("H129").Value = 2668.5588
$ws.Range("J129").Value = 2719.4243
$ws.Range("L129").Value = 8158.2729
$ws.Range("N129").Value = -18158.2729
# Row 132
$ws.Range("H132").Value = 4916.154
$ws.Range("I132").Value = 5562.222
$ws.Range("J132").Value = 3462.5
$ws.Range("K132").Value = 16686.666
$ws.Range("L132").Value = 10387.5
$ws.Range("M132").Value = -14156.666
$ws.Range("N132").Value = -15447.5
# Row 138
$ws.Range("H138").Value = 2159.2576
$ws.Range("I138").Value = 1520.2084
$ws.Range("J138").Value = 2524.4285
$ws.Range("K138").Value = 4560.6252
$ws.Range("L138").Value = 7573.2855
$ws.Range("M138").Value = 579.3747999999996
$ws.Range("N138").Value = -17853.2855
# Row 141
$ws.Range("H141").Value = 3812.125
$ws.Range("I141").Value = 2036
$ws.Range("J141").Value = 4619.4546
$ws.Range("K141").Value = 6108
$ws.Range("L141").Value = 13858.3638
$ws.Range("M141").Value = -928
$ws.Range("N141").Value = -24218.3638

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8710.391
$ws.Range("I32").Value = 8854.509
$ws.Range("J32").Value = 7536.857
$ws.Range("K32").Value = 8854.509
$ws.Range("L32").Value = 7536.857
$ws.Range("M32").Value = -8567.509
$ws.Range("N32").Value = -8110.857
# Row 88
$ws.Range("H88").Value = 2448.15
$ws.Range("I88").Value = 2446.6
$ws.Range("J88").Value = 2449.7
$ws.Range("K88").Value = 2446.6
$ws.Range("L88").Value = 2449.7
$ws.Range("M88").Value = -2040.6
$ws.Range("N88").Value = -3261.7
# Row 91
$ws.Range("H91").Value = 2448.15
$ws.Range("I91").Value = 2446.6
$ws.Range("J91").Value = 2449.7
$ws.Range("K91").Value = 2446.6
$ws.Range("L91").Value = 2449.7
$ws.Range("M91").Value = -1042.6
$ws.Range("N91").Value = -5257.7
# Row 122
$ws.Range("H122").Value = 13659.333
$ws.Range("I122").Value = 16833.715
$ws.Range("J122").Value = 2549
$ws.Range("K122").Value = 50501.145
$ws.Range("L122").Value = 7647
$ws.Range("M122").Value = -48051.145
$ws.Range("N122").Value = -12547

$ws = $wb.Worksheets.Item("BSM")
# Row 19
$ws.Range("H19").Value = 10000
$ws.Range("I19").Value = 10000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -9827
$ws.Range("N19").Value = ""
# Row 86
$ws.Range("H86").Value = 27780088
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 55558276
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 55558276
$ws.Range("M86").Value = -777
$ws.Range("N86").Value = -55560522
# Row 89
$ws.Range("H89").Value = 27780088
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 55558276
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 277791380
$ws.Range("M89").Value = -3884
$ws.Range("N89").Value = -277802612
# Row 107
$ws.Range("H107").Value = 2420.3333
$ws.Range("I107").Value = 2704.4
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2704.4
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -784.4000000000001
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
# Row 32
$ws.Range("H32").Value = 505377.5
$ws.Range("I32").Value = 505377.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 505377.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -505061.5
$ws.Range("N32").Value = ""
# Row 51
$ws.Range("H51").Value = 11099
$ws.Range("J51").Value = 11099
$ws.Range("L51").Value = 11099
$ws.Range("N51").Value = -12571
# Row 61
$ws.Range("H61").Value = 11099
$ws.Range("J61").Value = 11099
$ws.Range("L61").Value = 11099
$ws.Range("N61").Value = -11795
# Row 94
$ws.Range("H94").Value = 2550.7188
$ws.Range("J94").Value = 3800.6875
$ws.Range("L94").Value = 3800.6875
$ws.Range("N94").Value = -4702.6875
# Row 106
$ws.Range("H106").Value = 50085.5
$ws.Range("J106").Value = 50085.5
$ws.Range("L106").Value = 50085.5
$ws.Range("N106").Value = -52609.5
# Row 132
$ws.Range("H132").Value = 23812798
$ws.Range("I132").Value = 35716840
$ws.Range("K132").Value = 107150520
$ws.Range("M132").Value = -107147990

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 816.91
$ws.Range("I131").Value = 423.33334
$ws.Range("J131").Value = 855.83514
$ws.Range("K131").Value = 1270.00002
$ws.Range("L131").Value = 2567.50542
$ws.Range("M131").Value = 3769.99998
$ws.Range("N131").Value = -12647.50542
# Row 134
$ws.Range("H134").Value = 3652.353
$ws.Range("I134").Value = 1953.6364
$ws.Range("J134").Value = 6766.6665
$ws.Range("K134").Value = 5860.9092
$ws.Range("L134").Value = 20299.9995
$ws.Range("M134").Value = -790.9092000000001
$ws.Range("N134").Value = -30439.9995
# Row 138
$ws.Range("H138").Value = 7235.5
$ws.Range("I138").Value = 4530
$ws.Range("J138").Value = 7536.1113
$ws.Range("K138").Value = 13590
$ws.Range("L138").Value = 22608.3339
$ws.Range("M138").Value = -8450
$ws.Range("N138").Value = -32888.3339

$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
# Row 102
$ws.Range("H102").Value = 2622.862
$ws.Range("I102").Value = 2880.9583
$ws.Range("J102").Value = 1384
$ws.Range("K102").Value = 2880.9583
$ws.Range("L102").Value = 1384
$ws.Range("M102").Value = -1258.9583
$ws.Range("N102").Value = -4628

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2752.2173
$ws.Range("I122").Value = 2600.1052
$ws.Range("J122").Value = 3474.75
$ws.Range("K122").Value = 7800.3156
$ws.Range("L122").Value = 10424.25
$ws.Range("M122").Value = -5350.3156
$ws.Range("N122").Value = -15324.25
# Row 126
$ws.Range("H126").Value = 8799.875
$ws.Range("I126").Value = 2678.8
$ws.Range("J126").Value = 19001.666
$ws.Range("K126").Value = 8036.400000000001
$ws.Range("L126").Value = 57004.99800000001
$ws.Range("M126").Value = -5566.400000000001
$ws.Range("N126").Value = -61944.99800000001
# Row 132
$ws.Range("H132").Value = 1603.4
$ws.Range("I132").Value = 882.28125
$ws.Range("J132").Value = 4487.875
$ws.Range("K132").Value = 2646.84375
$ws.Range("L132").Value = 13463.625
$ws.Range("M132").Value = -116.84375
$ws.Range("N132").Value = -18523.625

Write-Output "Applied all Ultima_Profits updates."
